$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the two new components below the existing data (rows 3 and 4)
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "3,5 mm powerjack"

$ws.Range("A4").Value = 1
$ws.Range("B4").Value = "2,5 mm powerjack"

# Column B needs to widen to fit the new, longer text (bestFit width change)
$ws.Columns.Item(2).ColumnWidth = 16.5

# Update the active selection to reflect where the user ended up working
$ws.Range("C17").Select() | Out-Null
